$wb = $excel.ActiveWorkbook

# This report tracks localization handoff/handback status per source file,
# one row per file, one sheet per target language. A new handback event has
# landed for the "36c51d5e-0920-4625-bbbc-31c9a2100bcb" file (row 7) on both
# the zh-cn and de-de sheets. The handback turned out to be stale (not
# built from the latest source version), so "Latest Target File",
# "Latest Handback File", "Latest Handback DateTime" and "Error Detail"
# get populated accordingly.

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/68d486b298a7df48602ebbd1350752dff3ec0835/e2e/36c51d5e-0920-4625-bbbc-31c9a2100bcb.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/377d89d164c60fc4b44426c1d478e5f51cf3fbdf/e2e/36c51d5e-0920-4625-bbbc-31c9a2100bcb.md."
$hyperlinkTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/377d89d164c60fc4b44426c1d478e5f51cf3fbdf/e2e/36c51d5e-0920-4625-bbbc-31c9a2100bcb.md"
$targetFile = "36c51d5e-0920-4625-bbbc-31c9a2100bcb.md"

# ---------------------------------------------------------------------
# zh-cn sheet, row 7
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$i7Zh = $wsZh.Range("I7")
$i7Zh.Value = $targetFile
$i7Zh.Style = "HyperLink"
$wsZh.Hyperlinks.Add($i7Zh, $hyperlinkTarget, "", "", $targetFile)

$wsZh.Range("J7").Value = "36c51d5e-0920-4625-bbbc-31c9a2100bcb.c0a60e5d763b651e3f8c535753a4c24690f5e9a1.zh-cn.xlf"
$wsZh.Range("K7").Value = "2016-08-31 19:03:34"
$wsZh.Range("P7").Value = $errorDetail

# ---------------------------------------------------------------------
# de-de sheet, row 7
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$i7De = $wsDe.Range("I7")
$i7De.Value = $targetFile
$i7De.Style = "HyperLink"
$wsDe.Hyperlinks.Add($i7De, $hyperlinkTarget, "", "", $targetFile)

$wsDe.Range("J7").Value = "36c51d5e-0920-4625-bbbc-31c9a2100bcb.c0a60e5d763b651e3f8c535753a4c24690f5e9a1.de-de.xlf"
$wsDe.Range("K7").Value = "2016-08-31 19:03:41"
$wsDe.Range("P7").Value = $errorDetail

Write-Host "Handback report updated for row 7 on zh-cn and de-de sheets."
